# Commit: "Cleared all errors, also 25 students and 10 answer keys are
# displaying correctly"
#
# Effective changes in the workbook:
#   1. The "Answers " sheet (trailing space) is renamed to "Answers".
#   2. The "Answers" sheet becomes the active/selected tab (instead of
#      "Students"), which also clears the tabSelected flag that used to
#      sit on "Students".

$wb = $excel.ActiveWorkbook

$answers = $wb.Worksheets.Item(2)

# 1. Rename "Answers " -> "Answers"
$answers.Name = "Answers"

# 2. Make "Answers" the active sheet/tab (was "Students"). This also
#    clears tabSelected on "Students" as a side effect, since only one
#    sheet can be the selected tab at a time.
$answers.Activate()
$answers.Select()
